# Updates cryptos list values per Fri Nov 24 06:20:40 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.439.98"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.071.30"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'235.01"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'57.27"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value = "'0.392"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "2.375.94"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'14.42"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "'20.85"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "'0.779"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "2.072.84"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "37.389.54"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "'69.62"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "'227.16"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'167.84"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "'8.86"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "'1.40"
$ws.Range("E28").Value = "  -6.34%  "
$ws.Range("D29").Value = "'0.129"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").Value = "'19.08"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "'3.39"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.486.89"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0955"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").Value = "'96.95"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "'1.16"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "'4.11"
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "'15.24"
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("D49").Value = "'7.21"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.262.77"
$ws.Range("E51").Value = "  +0.38%  "
